$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.08244627342805799
$ws.Range("D2").Value = 0.04072730343476749
$ws.Range("E2").Value = 0.114761846639194
$ws.Range("F2").Value = 4.212240608295559
$ws.Range("G2").Value = 3.614556738805504
$ws.Range("H2").Value = 2.59203304739583
$ws.Range("I2").Value = 3.070072893264268
$ws.Range("J2").Value = 0.2435914505153747
$ws.Range("L2").Value = 0.2099675411597417
$ws.Range("C3").Value = 0.08274999149469409
$ws.Range("D3").Value = 0.03994293178559616
$ws.Range("E3").Value = 0.115101349042444
$ws.Range("F3").Value = 4.132352122186887
$ws.Range("G3").Value = 3.517155678698856
$ws.Range("H3").Value = 2.555584062300568
$ws.Range("I3").Value = 3.009488150405531
$ws.Range("J3").Value = 0.2427463245263866
$ws.Range("L3").Value = 0.2101991750564309
$ws.Range("C4").Value = 0.08295171857556305
$ws.Range("D4").Value = 0.03945342185512857
$ws.Range("E4").Value = 0.115350064938136
$ws.Range("F4").Value = 4.085917180444113
$ws.Range("G4").Value = 3.459749556916165
$ws.Range("H4").Value = 2.534735381654968
$ws.Range("I4").Value = 2.97421160859524
$ws.Range("J4").Value = 0.2423468978062289
$ws.Range("L4").Value = 0.2104273127015617
$ws.Range("C5").Value = 0.08303776512078898
$ws.Range("D5").Value = 0.03925191464291444
$ws.Range("E5").Value = 0.115461545523166
$ws.Range("F5").Value = 4.067648321696367
$ws.Range("G5").Value = 3.436954209929127
$ws.Range("H5").Value = 2.526622299106464
$ws.Range("I5").Value = 2.960316374123423
$ws.Range("J5").Value = 0.2422140927898866
$ws.Range("L5").Value = 0.2105418663030818
$ws.Range("C6").Value = 0.08305228532055686
$ws.Range("D6").Value = 0.0392183308581302
$ws.Range("E6").Value = 0.1154806685284235
$ws.Range("F6").Value = 4.064654153318969
$ws.Range("G6").Value = 3.433205039764971
$ws.Range("H6").Value = 2.525298195254777
$ws.Range("I6").Value = 2.958037992923124
$ws.Range("J6").Value = 0.2421938481046482
$ws.Range("L6").Value = 0.2105621910663373
$ws.Range("C7").Value = 0.08295286346701225
$ws.Range("D7").Value = 0.03945071252308452
$ws.Range("E7").Value = 0.1153515273962284
$ws.Range("F7").Value = 4.085668158771171
$ws.Range("G7").Value = 3.459439715698323
$ws.Range("H7").Value = 2.53462441806343
$ws.Range("I7").Value = 2.974022272345707
$ws.Range("J7").Value = 0.242344985529158
$ws.Range("L7").Value = 0.2104287702335732
$ws.Range("C8").Value = 0.08254783659716836
$ws.Range("D8").Value = 0.04045846743951031
$ws.Range("E8").Value = 0.1148705526879716
$ws.Range("F8").Value = 4.184149277594315
$ws.Range("G8").Value = 3.580471956284896
$ws.Range("H8").Value = 2.579146356522415
$ws.Range("I8").Value = 3.048782211994393
$ws.Range("J8").Value = 0.2432752020112616
$ws.Range("L8").Value = 0.210029565539557
$ws.Range("C9").Value = 0.08187417311723166
$ws.Range("D9").Value = 0.04237388412350995
$ws.Range("E9").Value = 0.1142467577313031
$ws.Range("F9").Value = 4.398246425471143
$ws.Range("G9").Value = 3.837093253837452
$ws.Range("H9").Value = 2.678706758664305
$ws.Range("I9").Value = 3.210802522471198
$ws.Range("J9").Value = 0.2460515404228687
$ws.Range("L9").Value = 0.209929476175013
$ws.Range("C10").Value = 0.08145228262232607
$ws.Range("D10").Value = 0.04374684716355404
$ws.Range("E10").Value = 0.113983230596947
$ws.Range("F10").Value = 4.568651750589652
$ws.Range("G10").Value = 4.037755198530874
$ws.Range("H10").Value = 2.759481704697123
$ws.Range("I10").Value = 3.339479224269638
$ws.Range("J10").Value = 0.2486782488693251
$ws.Range("L10").Value = 0.2102740097251683
$ws.Range("C11").Value = 0.08127612052316913
$ws.Range("D11").Value = 0.04436460860065239
$ws.Range("E11").Value = 0.1139056729406569
$ws.Range("F11").Value = 4.649092480800277
$ws.Range("G11").Value = 4.131757295952866
$ws.Range("H11").Value = 2.797920588989996
$ws.Range("I11").Value = 3.400165733711333
$ws.Range("J11").Value = 0.2500021465127134
$ws.Range("L11").Value = 0.2105219793650548
$ws.Range("C12").Value = 0.0812116711778863
$ws.Range("D12").Value = 0.04459760866350493
$ws.Range("E12").Value = 0.1138823922789776
$ws.Range("F12").Value = 4.679979194737001
$ws.Range("G12").Value = 4.167751078159426
$ws.Range("H12").Value = 2.812722783112577
$ws.Range("I12").Value = 3.423459690805061
$ws.Range("J12").Value = 0.2505221380716804
$ws.Range("L12").Value = 0.2106290343934134
$ws.Range("C13").Value = 0.08122545111196011
$ws.Range("D13").Value = 0.04454746867695292
$ws.Range("E13").Value = 0.1138871353343625
$ws.Range("F13").Value = 4.673308175592751
$ws.Range("G13").Value = 4.159981407583189
$ws.Range("H13").Value = 2.809523872179454
$ws.Range("I13").Value = 3.418428921644534
$ws.Range("J13").Value = 0.2504093168565475
$ws.Range("L13").Value = 0.2106053925704572
$ws.Range("C14").Value = 0.08127077299139174
$ws.Range("D14").Value = 0.04438379608434317
$ws.Range("E14").Value = 0.1139036355800833
$ws.Range("F14").Value = 4.651624989708097
$ws.Range("G14").Value = 4.134710529563506
$ws.Range("H14").Value = 2.799133424391812
$ws.Range("I14").Value = 3.402075839405256
$ws.Range("J14").Value = 0.2500445519021923
$ws.Range("L14").Value = 0.2105305230119541
$ws.Range("C15").Value = 0.08129882800266053
$ws.Range("D15").Value = 0.0442834217970649
$ws.Range("E15").Value = 0.1139145354882327
$ws.Range("F15").Value = 4.638398991914102
$ws.Range("G15").Value = 4.119283301500786
$ws.Range("H15").Value = 2.792801120007255
$ws.Range("I15").Value = 3.392100025023069
$ws.Range("J15").Value = 0.2498235563125561
$ws.Range("L15").Value = 0.2104863773735062
$ws.Range("C16").Value = 0.0814641117557624
$ws.Range("D16").Value = 0.04370634195043266
$ws.Range("E16").Value = 0.1139891510305127
$ws.Range("F16").Value = 4.563454017970486
$ws.Range("G16").Value = 4.03166719115859
$ws.Range("H16").Value = 2.757003943429652
$ws.Range("I16").Value = 3.335556838464754
$ws.Range("J16").Value = 0.2485943341445704
$ws.Range("L16").Value = 0.2102596431740693
$ws.Range("C17").Value = 0.08156953925634802
$ws.Range("D17").Value = 0.0433506113365496
$ws.Range("E17").Value = 0.1140457667778136
$ws.Range("F17").Value = 4.518229954706555
$ws.Range("G17").Value = 3.978618488073664
$ws.Range("H17").Value = 2.7354792489337
$ws.Range("I17").Value = 3.301423122165716
$ws.Range("J17").Value = 0.2478733625171117
$ws.Range("L17").Value = 0.21014394140159
$ws.Range("C18").Value = 0.08163166196550797
$ws.Range("D18").Value = 0.0431453610546626
$ws.Range("E18").Value = 0.1140823140050493
$ws.Range("F18").Value = 4.49249302299387
$ws.Range("G18").Value = 3.948361920073637
$ws.Range("H18").Value = 2.723258095245626
$ws.Range("I18").Value = 3.281992512070133
$ws.Range("J18").Value = 0.2474708134322796
$ws.Range("L18").Value = 0.2100859785397517
$ws.Range("C19").Value = 0.08165295066009293
$ws.Range("D19").Value = 0.04307575493768567
$ws.Range("E19").Value = 0.1140953723361378
$ws.Range("F19").Value = 4.483825969884379
$ws.Range("G19").Value = 3.938161289513062
$ws.Range("H19").Value = 2.719147500196925
$ws.Range("I19").Value = 3.275448249123514
$ws.Range("J19").Value = 0.2473365980628799
$ws.Range("L19").Value = 0.2100678269108585
$ws.Range("C20").Value = 0.0815581628205706
$ws.Range("D20").Value = 0.04338854573108009
$ws.Range("E20").Value = 0.1140393276552434
$ws.Range("F20").Value = 4.523015659545337
$ws.Range("G20").Value = 3.984239111306863
$ws.Range("H20").Value = 2.737754085559743
$ws.Range("I20").Value = 3.305035758850011
$ws.Range("J20").Value = 0.2479488544886479
$ws.Range("L20").Value = 0.2101553692208284
$ws.Range("C21").Value = 0.08125739958887657
$ws.Range("D21").Value = 0.04443189563728112
$ws.Range("E21").Value = 0.1138986237849053
$ws.Range("F21").Value = 4.657982274591717
$ws.Range("G21").Value = 4.142122368998002
$ws.Range("H21").Value = 2.802178646993468
$ws.Range("I21").Value = 3.406870600930233
$ws.Range("J21").Value = 0.2501511848732676
$ws.Range("L21").Value = 0.210552156729598
$ws.Range("C22").Value = 0.08107400013837207
$ws.Range("D22").Value = 0.04510837003721591
$ws.Range("E22").Value = 0.1138421573596506
$ws.Range("F22").Value = 4.748673783272181
$ws.Range("G22").Value = 4.247626492829795
$ws.Range("H22").Value = 2.845720272745666
$ws.Range("I22").Value = 3.475253603156972
$ws.Range("J22").Value = 0.2516993557020299
$ws.Range("L22").Value = 0.2108881709999224
$ws.Range("C23").Value = 0.08117068124059301
$ws.Range("D23").Value = 0.04474780300201076
$ws.Range("E23").Value = 0.1138690459237601
$ws.Range("F23").Value = 4.700041032170901
$ws.Range("G23").Value = 4.191102778043899
$ws.Range("H23").Value = 2.822348962675051
$ws.Range("I23").Value = 3.438587690463493
$ws.Range("J23").Value = 0.2508630729908106
$ws.Range("L23").Value = 0.2107018048311602
$ws.Range("C24").Value = 0.08156330139994239
$ws.Range("D24").Value = 0.04337139788746924
$ws.Range("E24").Value = 0.1140422263295857
$ws.Range("F24").Value = 4.520851223648407
$ws.Range("G24").Value = 3.981697274730891
$ws.Range("H24").Value = 2.73672515344083
$ws.Range("I24").Value = 3.303401883963858
$ws.Range("J24").Value = 0.2479146873642009
$ws.Range("L24").Value = 0.2101501760574394
$ws.Range("C25").Value = 0.0820435555268002
$ws.Range("D25").Value = 0.04186193773043811
$ws.Range("E25").Value = 0.1143813132760521
$ws.Range("F25").Value = 4.338049164424831
$ws.Range("G25").Value = 3.765568830016434
$ws.Range("H25").Value = 2.650445154802867
$ws.Range("I25").Value = 3.165296621382623
$ws.Range("J25").Value = 0.2451978751497634
$ws.Range("L25").Value = 0.2098832758478721
